$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 616 (shifts existing rows 616.. down by one)
$ws.Rows.Item(616).Insert()

# Populate the newly inserted row 616 with the new data point
$ws.Cells.Item(616, 1).NumberFormat = "@"
$ws.Cells.Item(616, 1).Value = "2026/01/13"
$ws.Cells.Item(616, 1).Style = "Normal"
$ws.Cells.Item(616, 2).Value = "火"
$ws.Cells.Item(616, 3).Value = 17
$ws.Cells.Item(616, 4).Value = 27
